$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-6
# from serial 45221 (2023-10-22) to serial 45224 (2023-10-25)
$ws.Range("C2:C6").Value = 45224
